# Applies the "5.0.0 -> 6.0.0" release bump to the Span StructureDefinition
# workbook: Version, Date, Publisher/Jurisdiction metadata, and the
# root-element Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (previously empty)
$meta.Range("B9").Value = "Alvearie Team"

# The sheet had a duplicated "Contact" row (10 and 11 were identical); drop
# one of them so every row below shifts up by one without us having to
# retype their (unchanged) content/styling.
$meta.Rows.Item(10).Delete()

# The remaining former "Contact" row becomes the new "Jurisdiction" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# The root "Extension" element's Short/Definition columns (K/L) now carry the
# profile-specific text instead of the generic placeholders.
$elements.Range("K2").Value = "Span"
$elements.Range("L2").Value = "Detail on a span of text from a reference source used as input for an insight evaluation"
